# Fixed all errors found by manually reviewing the master BOM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 - Cap, 0402-A, 12pF: corrected manufacturer part number and price
$ws.Range("C7").Value = "CC0402JRNPO9BN120"
$ws.Range("G7").Value = 0.0135
$ws.Range("H7").Value = 0.135

# Row 11 - Cap, 6-0805_M, 10uF: corrected manufacturer part number and price
$ws.Range("C11").Value = "GRM21BR61E106KA73K"
$ws.Range("G11").Value = 0.1635
$ws.Range("H11").Value = 0.327

# Row 14 - Cap, 0402-A, 47nF: corrected manufacturer part number (was a duplicate
# of row 7's old, wrong part) and price
$ws.Range("C14").Value = "GRM155R71H473KE14D"
$ws.Range("G14").Value = 0.036
$ws.Range("H14").Value = 0.036

# Row 18 - Cap, 0402-A, 150nF: corrected manufacturer part number and price
$ws.Range("C18").Value = "CL05A154KP5NNNC"
$ws.Range("G18").Value = 0.0135
$ws.Range("H18").Value = 0.0135

# Row 23 - B4B-XH-A connector: corrected price only
$ws.Range("G23").Value = 0.21
$ws.Range("H23").Value = 0.21

# Row 34 - Res1, 6-0805_M, 330: corrected manufacturer part number and price
$ws.Range("C34").Value = "CRGCQ0805F330R"
$ws.Range("G34").Value = 0.0333
$ws.Range("H34").Value = 0.0333

# Row 36 - Res1, 0402-A, 150k: corrected manufacturer part number and price
$ws.Range("C36").Value = "PFR05S-154-FNH"
$ws.Range("G36").Value = 0.1
$ws.Range("H36").Value = 0.1

# Row 38 - Res1, 0402-A, 78.7K -> 78.7k: corrected manufacturer part number,
# value casing, and price
$ws.Range("C38").Value = "RC0402FR-0778K7L"
$ws.Range("F38").Value = "78.7k"
$ws.Range("G38").Value = 0.0959
$ws.Range("H38").Value = 0.0959
